$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Top block (rows 3-7): standard errors for j = 1..5
$ws.Range("B3").Value = "(0.244)"
$ws.Range("C3").Value = "(0.251)"

$ws.Range("B4").Value = "(0.268)"
$ws.Range("C4").Value = "(0.281)"

$ws.Range("B5").Value = "(0.310)"
$ws.Range("C5").Value = "(0.326)"

$ws.Range("B6").Value = "(0.362)"
$ws.Range("C6").Value = "(0.379)"

$ws.Range("B7").Value = "(0.417)"
$ws.Range("C7").Value = "(0.435)"

# x / j=1 block (row 9)
$ws.Range("B9").Value = "(0.287)"
$ws.Range("C9").Value = "(0.287)"

# var(M1[i]) row 13
$ws.Range("B13").Value = "(0.636)"
$ws.Range("C13").Value = "(0.638)"

# var(M2[i>id]) row 14
$ws.Range("B14").Value = "(2.161)"
$ws.Range("C14").Value = "(2.164)"

# var(e.yobs) row 15
$ws.Range("B15").Value = "(0.957)"
$ws.Range("C15").Value = "(0.957)"

# second x block (rows 17-21), column C only
$ws.Range("C17").Value = "(0.058)"
$ws.Range("C18").Value = "(0.030)"
$ws.Range("C19").Value = "(0.005)"
$ws.Range("C20").Value = "(0.035)"
$ws.Range("C21").Value = "(0.025)"

# Number of observations row 22
$ws.Range("B22").Value = 8480
$ws.Range("C22").Value = 10352
